$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.423576
$ws.Range("H2").Value = 19.270728
$ws.Range("I2").Value = 0.001681024218962088
$ws.Range("J2").Value = 0.001681024218962088
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 635.591801565024
$ws.Range("R2").Value = 5720.326214085216
$ws.Range("S2").Value = 0.0003526805105533398
$ws.Range("T2").Value = 0.0003526805105533399
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.423576
$ws.Range("H3").Value = 19.270728
$ws.Range("I3").Value = 0.001681024218962088
$ws.Range("J3").Value = 0.001681024218962088
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 1047.082943278744
$ws.Range("R3").Value = 9423.746489508696
$ws.Range("S3").Value = 0.000581010872257863
$ws.Range("T3").Value = 0.0005810108722578631
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.423576
$ws.Range("H4").Value = 19.270728
$ws.Range("I4").Value = 0.001681024218962088
$ws.Range("J4").Value = 0.001681024218962088
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 420.0768612774401
$ws.Range("R4").Value = 3780.69175149696
$ws.Range("S4").Value = 0.0002330944507814192
$ws.Range("T4").Value = 0.0002330944507814193
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.423576
$ws.Range("H5").Value = 19.270728
$ws.Range("I5").Value = 0.001681024218962088
$ws.Range("J5").Value = 0.001681024218962088
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 926.74727412088
$ws.Range("R5").Value = 8340.725467087919
$ws.Range("S5").Value = 0.0005142383853694659
$ws.Range("T5").Value = 0.000514238385369466
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3580.644531333333
$ws.Range("H6").Value = 10741.933594
$ws.Range("I6").Value = 0.9370403925578976
$ws.Range("J6").Value = 0.9370403925578976
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 354293.0461839487
$ws.Range("R6").Value = 3188637.415655538
$ws.Range("S6").Value = 0.196591982630962
$ws.Range("T6").Value = 0.196591982630962
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3580.644531333333
$ws.Range("H7").Value = 10741.933594
$ws.Range("I7").Value = 0.9370403925578976
$ws.Range("J7").Value = 0.9370403925578976
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 583667.3863130825
$ws.Range("R7").Value = 5253006.476817743
$ws.Range("S7").Value = 0.3238684188363814
$ws.Range("T7").Value = 0.3238684188363815
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3580.644531333333
$ws.Range("H8").Value = 10741.933594
$ws.Range("I8").Value = 0.9370403925578976
$ws.Range("J8").Value = 0.9370403925578976
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 234160.2117064914
$ws.Range("R8").Value = 2107441.905358423
$ws.Range("S8").Value = 0.1299320457132656
$ws.Range("T8").Value = 0.1299320457132656
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3580.644531333333
$ws.Range("H9").Value = 10741.933594
$ws.Range("I9").Value = 0.9370403925578976
$ws.Range("J9").Value = 0.9370403925578976
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 516589.600404666
$ws.Range("R9").Value = 4649306.403641994
$ws.Range("S9").Value = 0.2866479453772885
$ws.Range("T9").Value = 0.2866479453772886
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.9157713333333334
$ws.Range("H10").Value = 2.747314
$ws.Range("I10").Value = 0.0002396537054071653
$ws.Range("J10").Value = 0.0002396537054071653
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 90.61257336644535
$ws.Range("R10").Value = 815.5131602980081
$ws.Range("S10").Value = 0.00005027957969052017
$ws.Range("T10").Value = 0.00005027957969052018
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.9157713333333334
$ws.Range("H11").Value = 2.747314
$ws.Range("I11").Value = 0.0002396537054071653
$ws.Range("J11").Value = 0.0002396537054071653
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 149.2764377781109
$ws.Range("R11").Value = 1343.487940002998
$ws.Range("S11").Value = 0.00008283129228466298
$ws.Range("T11").Value = 0.00008283129228466301
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.9157713333333334
$ws.Range("H12").Value = 2.747314
$ws.Range("I12").Value = 0.0002396537054071653
$ws.Range("J12").Value = 0.0002396537054071653
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 59.88787979694224
$ws.Range("R12").Value = 538.9909181724802
$ws.Range("S12").Value = 0.00003323090066727651
$ws.Range("T12").Value = 0.00003323090066727652
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.9157713333333334
$ws.Range("H13").Value = 2.747314
$ws.Range("I13").Value = 0.0002396537054071653
$ws.Range("J13").Value = 0.0002396537054071653
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 132.1208913671622
$ws.Range("R13").Value = 1189.08802230446
$ws.Range("S13").Value = 0.00007331193276470557
$ws.Range("T13").Value = 0.00007331193276470559
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 233.243637
$ws.Range("H14").Value = 699.7309110000001
$ws.Range("I14").Value = 0.0610389295177331
$ws.Range("J14").Value = 0.06103892951773311
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 23078.69377499519
$ws.Range("R14").Value = 207708.2439749567
$ws.Range("S14").Value = 0.01280602657779372
$ws.Range("T14").Value = 0.01280602657779372
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 233.243637
$ws.Range("H15").Value = 699.7309110000001
$ws.Range("I15").Value = 0.0610389295177331
$ws.Range("J15").Value = 0.06103892951773311
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 38020.16726057246
$ws.Range("R15").Value = 342181.5053451521
$ws.Range("S15").Value = 0.02109682970699909
$ws.Range("T15").Value = 0.0210968297069991
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 233.243637
$ws.Range("H16").Value = 699.7309110000001
$ws.Range("I16").Value = 0.0610389295177331
$ws.Range("J16").Value = 0.06103892951773311
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 15253.22576457329
$ws.Range("R16").Value = 137279.0318811595
$ws.Range("S16").Value = 0.008463789867945165
$ws.Range("T16").Value = 0.008463789867945166
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 233.243637
$ws.Range("H17").Value = 699.7309110000001
$ws.Range("I17").Value = 0.0610389295177331
$ws.Range("J17").Value = 0.06103892951773311
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 33650.71181469481
$ws.Range("R17").Value = 302856.4063322534
$ws.Range("S17").Value = 0.01867228336499511
$ws.Range("T17").Value = 0.01867228336499512
